$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Sandeep Sharma"

# Insert a new column before column A for "matchNo" (shifts B..L -> C..M)
$ws.Range("A1").EntireColumn.Insert()

# Header row
$ws.Range("A1").Value = "matchNo"
$ws.Range("B1").Value = "teamName"
$ws.Range("C1").Value = "batterName"
$ws.Range("D1").Value = "states"
$ws.Range("E1").Value = "runs"
$ws.Range("F1").Value = "balls"
$ws.Range("G1").Value = "fours"
$ws.Range("H1").Value = "sixes"
$ws.Range("I1").Value = "sr"
$ws.Range("J1").Value = "opponentTeamName"
$ws.Range("K1").Value = "venue"
$ws.Range("L1").Value = "date"
$ws.Range("M1").Value = "result"

# Row 2 - add matchNo value (rest already shifted via insert); force text
# storage ("'" prefix) for values Excel would otherwise coerce to numbers.
$ws.Range("A2").Value = "28th"
$ws.Range("B2").Value = "Sunrisers Hyderabad"
$ws.Range("C2").Value = "Sandeep Sharma"
$ws.Range("D2").Value = "'"
$ws.Range("E2").Value = "'8"
$ws.Range("F2").Value = "'6"
$ws.Range("G2").Value = "'1"
$ws.Range("H2").Value = "'0"
$ws.Range("I2").Value = "'133.33"
$ws.Range("J2").Value = "Rajasthan Royals"
$ws.Range("K2").Value = "Delhi"
$ws.Range("L2").Value = "May 02"
$ws.Range("M2").Value = "Royals won by 55 runs"

# Row 3 - new row
$ws.Range("A3").Value = "33rd"
$ws.Range("B3").Value = "Sunrisers Hyderabad"
$ws.Range("C3").Value = "Sandeep Sharma"
$ws.Range("D3").Value = "run out (†Pant)"
$ws.Range("E3").Value = "'0"
$ws.Range("F3").Value = "'1"
$ws.Range("G3").Value = "'0"
$ws.Range("H3").Value = "'0"
$ws.Range("I3").Value = "'0.00"
$ws.Range("J3").Value = "Delhi Capitals"
$ws.Range("K3").Value = "Dubai (DSC)"
$ws.Range("L3").Value = "September 22"
$ws.Range("M3").Value = "Capitals won by 8 wickets (with 13 balls remaining)"
